$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (England)
$ws.Range("B2").Value = 720
$ws.Range("C2").Value = 2311
$ws.Range("D2").Value = 76
$ws.Range("E2").Value = 83

# Row 3 (Northern Ireland)
$ws.Range("B3").Value = 38
$ws.Range("C3").Value = 60

# Row 4 (Scotland)
$ws.Range("B4").Value = 173
$ws.Range("C4").Value = 415
$ws.Range("D4").Value = 28
$ws.Range("E4").Value = 7

# Row 5 (Wales)
$ws.Range("B5").Value = 76
$ws.Range("C5").Value = 168
